$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $cellRef, $val) {
    # Force the cell to be written as a literal text string so that
    # numeric-looking values (e.g. "291.17") are not auto-converted
    # into actual numbers by Excel, matching the inlineStr cells in the
    # original workbook. Resetting the Style afterwards avoids leaving
    # a stray "@" number-format style applied to the cell.
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue $ws 'D2' '22.380.98'
Set-TextValue $ws 'E2' '  -4.42%  '
Set-TextValue $ws 'D3' '1.567.49'
Set-TextValue $ws 'E3' '  -4.75%  '
Set-TextValue $ws 'D4' '1.002'
Set-TextValue $ws 'E4' '  -0.20%  '
Set-TextValue $ws 'E5' '  -0.09%  '
Set-TextValue $ws 'D6' '291.17'
Set-TextValue $ws 'E6' '  -2.74%  '
Set-TextValue $ws 'D7' '0.3666'
Set-TextValue $ws 'E7' '  -3.31%  '
Set-TextValue $ws 'D8' '49.53'
Set-TextValue $ws 'E8' '  -1.02%  '
Set-TextValue $ws 'D9' '0.3374'
Set-TextValue $ws 'E9' '  -5.17%  '
Set-TextValue $ws 'D10' '1.168'
Set-TextValue $ws 'E10' '  -4.06%  '
Set-TextValue $ws 'D11' '0.07568'
Set-TextValue $ws 'E11' '  -6.54%  '
Set-TextValue $ws 'E12' '  -0.11%  '
Set-TextValue $ws 'E13' '  -4.25%  '
Set-TextValue $ws 'D14' '6.047'
Set-TextValue $ws 'E14' '  -5.40%  '
Set-TextValue $ws 'E15' '  -6.72%  '
Set-TextValue $ws 'D16' '0.00001143'
Set-TextValue $ws 'E16' '  -4.24%  '
Set-TextValue $ws 'D17' '1.571.11'
Set-TextValue $ws 'E17' '  -4.87%  '
Set-TextValue $ws 'D18' '89.09'
Set-TextValue $ws 'E18' '  -8.55%  '
Set-TextValue $ws 'D19' '0.06714'
Set-TextValue $ws 'E19' '  -3.44%  '
Set-TextValue $ws 'E20' '  +0.02%  '
Set-TextValue $ws 'D21' '6.249'
Set-TextValue $ws 'E21' '  -7.41%  '
Set-TextValue $ws 'B22' 'Avalanche'
Set-TextValue $ws 'C22' 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue $ws 'D22' '16.38'
Set-TextValue $ws 'E22' '  -5.04%  '
Set-TextValue $ws 'B23' 'BitDAO'
Set-TextValue $ws 'C23' 'https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit'
Set-TextValue $ws 'D23' '0.5268'
Set-TextValue $ws 'E23' '  -8.69%  '
Set-TextValue $ws 'D24' '11.97'
Set-TextValue $ws 'E24' '  -3.43%  '
Set-TextValue $ws 'D25' '22.408.12'
Set-TextValue $ws 'E25' '  -4.37%  '
Set-TextValue $ws 'D26' '2.382'
Set-TextValue $ws 'E26' '  -4.89%  '
Set-TextValue $ws 'D27' '2.922'
Set-TextValue $ws 'E27' '  +0.66%  '
Set-TextValue $ws 'D28' '19.89'
Set-TextValue $ws 'E28' '  -4.66%  '
Set-TextValue $ws 'D29' '146.14'
Set-TextValue $ws 'E29' '  -4.74%  '
Set-TextValue $ws 'D30' '4.952'
Set-TextValue $ws 'E30' '  -4.98%  '
Set-TextValue $ws 'D31' '124.85'
Set-TextValue $ws 'E31' '  -5.94%  '
Set-TextValue $ws 'D32' '1.749.43'
Set-TextValue $ws 'E32' '  -4.66%  '
Set-TextValue $ws 'D33' '6.239'
Set-TextValue $ws 'E33' '  -9.75%  '
Set-TextValue $ws 'D34' '1.979'
Set-TextValue $ws 'E34' '  -6.48%  '
Set-TextValue $ws 'D35' '0.9821'
Set-TextValue $ws 'E35' '  -3.02%  '
Set-TextValue $ws 'E36' '  -12.44%  '
Set-TextValue $ws 'D37' '0.08412'
Set-TextValue $ws 'E37' '  -3.69%  '
Set-TextValue $ws 'D38' '0.02527'
Set-TextValue $ws 'E38' '  -7.06%  '
Set-TextValue $ws 'D39' '0.2294'
Set-TextValue $ws 'E39' '  -5.58%  '
Set-TextValue $ws 'D40' '0.06499'
Set-TextValue $ws 'E40' '  -4.10%  '
Set-TextValue $ws 'D41' '5.508'
Set-TextValue $ws 'E41' '  -7.08%  '
Set-TextValue $ws 'E42' '  -10.01%  '
Set-TextValue $ws 'E43' '  -5.12%  '
Set-TextValue $ws 'D44' '0.6379'
Set-TextValue $ws 'E44' '  -7.06%  '
Set-TextValue $ws 'D45' '14.54'
Set-TextValue $ws 'E45' '  -6.08%  '
Set-TextValue $ws 'E46' '  -0.10%  '
Set-TextValue $ws 'D47' '0.6019'
Set-TextValue $ws 'E47' '  -5.63%  '
Set-TextValue $ws 'D48' '3.773'
Set-TextValue $ws 'E48' '  -3.68%  '
Set-TextValue $ws 'D49' '2.114'
Set-TextValue $ws 'E49' '  -5.93%  '
Set-TextValue $ws 'D50' '121.04'
Set-TextValue $ws 'E50' '  -4.89%  '
Set-TextValue $ws 'E51' '  -5.75%  '
